$d = $word.ActiveDocument

# --- remove the _GoBack bookmark from its original location (end of doc) ---
# (a fresh one is re-inserted at paragraph 7 below; deleting it first avoids
#  touching the final paragraph mark, which Word will not let us delete)
try {
    $d.Bookmarks.Item("_GoBack").Delete()
} catch {
    # no pre-existing _GoBack bookmark -- nothing to do
}

# --- paragraph 14: Applications of graph theory ---
$p = $d.Paragraphs.Item(14)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Applications of </w:t></w:r><w:r><w:t>graph theory</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:br/><w:t>-</w:t></w:r><w:r><w:tab/></w:r><w:r><w:t>t</w:t></w:r><w:r><w:t>ravelling sales man, quantum com</w:t></w:r><w:r><w:t>puting</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $p.Range.InsertXML($xml)

# --- paragraph 10: graph theory big paragraph edits ---
$p = $d.Paragraphs.Item(10)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1080" w:firstLine="360"/></w:pPr><w:r><w:t xml:space="preserve">Graph theory is a discipline within mathematics that handles the study of mathematical structures known as </w:t></w:r><w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>graphs</w:t></w:r><w:r><w:t xml:space="preserve">. A graph is defined as a series of related nodes, or </w:t></w:r><w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>vertices</w:t></w:r><w:r><w:t>, which represent entities or states within a given mathematical</w:t></w:r><w:r><w:t xml:space="preserve"> system. Relationships between</w:t></w:r><w:r><w:t xml:space="preserve"> vertices of a graph are represented by </w:t></w:r><w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>edges</w:t></w:r><w:r><w:t xml:space="preserve">, which are illustrated in a graph </w:t></w:r><w:r><w:t xml:space="preserve">using lines or curves between a </w:t></w:r><w:r><w:t>set of vertices.</w:t></w:r><w:r><w:t xml:space="preserve"> An edge which connects a vertex to itself is called a </w:t></w:r><w:r><w:rPr><w:b/><w:i/></w:rPr><w:t>loop</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> A vertex of a graph is said to have </w:t></w:r><w:r><w:rPr><w:b/><w:i/></w:rPr><w:t xml:space="preserve">nth degree </w:t></w:r><w:r><w:t>when said ve</w:t></w:r><w:r><w:t xml:space="preserve">rtex contains n-number of edges connected to it. </w:t></w:r><w:r><w:br/></w:r><w:r><w:lastRenderedPageBreak/><w:tab/><w:t>The implied goal of employing a graph-theoretic approach to a given problem i</w:t></w:r><w:r><w:t>s to represent a set of data, and</w:t></w:r><w:r><w:t xml:space="preserve"> any relationships inherent within that set</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> in a more intuitive way, and </w:t></w:r><w:r><w:t xml:space="preserve">to create a model of a problem in such a way that one may find more profound understanding of the features and characteristics governing the problem. </w:t></w:r><w:r><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $p.Range.InsertXML($xml)

# --- paragraph 9: delete '(introducing ...)' paragraph ---
$p = $d.Paragraphs.Item(9)
$null = $p.Range.Delete()

# --- paragraph 8: Overview of graph theory (unchanged content) ---
$p = $d.Paragraphs.Item(8)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Overview of graph theory: </w:t></w:r><w:r><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $p.Range.InsertXML($xml)

# --- paragraph 7: Instant Insanity puzzle and solution description** ---
$p = $d.Paragraphs.Item(7)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Instant Insanity puzzle and solution description**</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $p.Range.InsertXML($xml)

# --- paragraph 6: difficulty -> probability text + spacer ---
$p = $d.Paragraphs.Item(6)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Given the fact that there is a single solution for a configuration of cubes, the probability of generating a solution in a completely random way would be 1 in 41,472, or roughly a 0.002% chance. This is more than three times less likely than the chance you will be struck by lightning in your lifetime, which is quoted as about a 1 in 13,000 chance (as per the National Oceanic and Atmospheric Administration: </w:t></w:r><w:r><w:t>http://www.lightningsafety.noaa.gov/odds.s</w:t></w:r><w:r><w:t>h</w:t></w:r><w:r><w:t>tml</w:t></w:r><w:r><w:t>).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1440"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $p.Range.InsertXML($xml)

# --- paragraph 5: number of possibilities -> full explanation ---
$p = $d.Paragraphs.Item(5)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Given the problem includes 4 separate cubes and 4 unique colors, the total number of possible arrangements is </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>41,472</w:t></w:r><w:r><w:t>. This number is calculated as follows:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="2160"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>(1)</w:t></w:r><w:r><w:t xml:space="preserve"> There are 3 ways in which the first cube may be chosen, by first deciding </w:t></w:r><w:r><w:t>which pair of faces will be left-right pair. Only 3 pairs exist on a cube, therefore the first choice is simplest.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="2160"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>(2)</w:t></w:r><w:r><w:t xml:space="preserve"> This number is then multiplied by the number of ways in which a following cube may be added to the subsequent in relation to it. In this case, there are 24 ways in which a new cube may be added in relation to an established sequence of cubes. This is then repeated for each remaining cube to be added to the sequence, finally generating a total of </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>3 x 24 x 24 x 24 = 41,472.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="2160"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $p.Range.InsertXML($xml)

# --- paragraph 4: Summary of the problem (drop colon) + spacer ---
$p = $d.Paragraphs.Item(4)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Summary of the problem</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1440"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $p.Range.InsertXML($xml)

# --- paragraph 3: Introduction + br + new presentation paragraph ---
$p = $d.Paragraphs.Item(3)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Introduction</w:t></w:r><w:r><w:br/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="1440" w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">In this presentation, we will apply graph theoretic solutions to a notable combinatorial problem involving a set of four cubes with colored faces. The problem was condensed into a game and sold under a few different names in North American market in th</w:t></w:r><w:r><w:t>e past century. Most popularly dubbed</w:t></w:r><w:r><w:t xml:space="preserve"> “Instant Insanity” or “Th</w:t></w:r><w:r><w:t>e Great Tantalizer”, there exist</w:t></w:r><w:r><w:t xml:space="preserve"> varying versions of the game, </w:t></w:r><w:r><w:t>however most follow</w:t></w:r><w:r><w:t xml:space="preserve"> suit of the original configuration, whic</w:t></w:r><w:r><w:t>h allowed for only one solution.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>S</w:t></w:r><w:r><w:t>ome</w:t></w:r><w:r><w:t xml:space="preserve"> other versions, although,</w:t></w:r><w:r><w:t xml:space="preserve"> contain</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>their own unique</w:t></w:r><w:r><w:t xml:space="preserve"> color</w:t></w:r><w:r><w:t xml:space="preserve"> configurations and number of solutions to the general problem. </w:t></w:r><w:r><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $p.Range.InsertXML($xml)

# --- paragraph 2: merge author-line runs ---
$p = $d.Paragraphs.Item(2)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Nicholas Brandt, Yi Yi (Lily) Zhang, Thoa Ta</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $p.Range.InsertXML($xml)
